$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix existing row 102: correct the date/time value in column A ---
$ws.Range("A102").Value = 45483.2916666667

# --- Append new row 103 with the latest day's OHLCV data ---

# Column A: date (copy number-format/style from A102 so style index is reused, not duplicated)
$ws.Range("A102").Copy()
$ws.Range("A103").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A103").Value = 45484.6436689815

# Columns B-F: plain numeric values
$ws.Range("B103").Value = 1800
$ws.Range("C103").Value = 6
$ws.Range("D103").Value = 5.96000003814697
$ws.Range("E103").Value = 6
$ws.Range("F103").Value = 6

# Column G: adj_close is stored as text (matches existing column G cells, e.g. G102)
$ws.Range("G103").NumberFormat = "@"
$ws.Range("G103").Value = "6"
$ws.Range("G103").ClearFormats()

# Column H: ticker text
$ws.Range("H103").Value = "PAL.MI"
